$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price/hour strings remain stored as text (matching original inline-string cells)
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# --- Column D (Price) updates ---
$ws.Range("D2").Value = "249.62"
$ws.Range("D4").Value = "5.429"
$ws.Range("D5").Value = "0.05624"
$ws.Range("D6").Value = "3.424"
$ws.Range("D7").Value = "6.367"
$ws.Range("D8").Value = "0.8115"
$ws.Range("D9").Value = "0.8934"
$ws.Range("D11").Value = "0.07499"
$ws.Range("D12").Value = "0.03087"
$ws.Range("D13").Value = "0.03087"
$ws.Range("D14").Value = "0.09321"
$ws.Range("D15").Value = "3.557"
$ws.Range("D16").Value = "0.001606"
$ws.Range("D17").Value = "0.04740"
$ws.Range("D18").Value = "0.0005795"
$ws.Range("D19").Value = "0.006411"
$ws.Range("D20").Value = "0.004993"
$ws.Range("D21").Value = "0.001031"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("D24").Value = "2.177"
$ws.Range("D25").Value = "0.3301"
$ws.Range("D28").Value = "0.0003008"
$ws.Range("D40").Value = "0.04040"
$ws.Range("D41").Value = "0.006836"
$ws.Range("D43").Value = "0.002724"
$ws.Range("D44").Value = "0.007511"
$ws.Range("D45").Value = "0.00005576"
$ws.Range("D47").Value = "0.5004"
$ws.Range("D48").Value = "0.2401"
$ws.Range("D49").Value = "0.00002102"

# --- Column E (Volume(1h)) updates: "Worstin24h" suffix moved from One(row18) to AAXToken(row27) ---
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E27").Value = "26AAXTokenAABWorstin24h"

# --- Column G (Hora) updates: hour changed from 17 to 18 for all data rows ---
$ws.Range("G2:G51").Value = "18"

# Restore default style/number format on touched ranges (remove temporary text-format marker)
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
